# "inclusao dos valores de tensao maxima no excel"
# Adds the peak-voltage / RMS-voltage / power readings for the extra
# measurement point (rows 11-13) below the existing measurement table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: peak voltage reading (V pico), next to the existing
#     "Ponto de tensao maxima na carga" / "6,2 cm" label pair in A11:B11.
$ws.Range("C11").Value = 7.3

# Labels are entered D13, D12, D11 (in that order) so the new shared
# strings land at the same table indices as the source workbook
# (Watts=11, V RMS=12, V pico=13).
$ws.Range("D13").Value = "Watts"
$ws.Range("D12").Value = "V RMS"
$ws.Range("D11").Value = "V pico"

# --- Row 12: RMS voltage derived from the peak voltage in C11.
$ws.Range("C12").Formula = "=C11/SQRT(2)"
$ws.Range("C12").NumberFormat = "0.000"

# --- Row 13: power delivered to the load (Ohms in H3) at that RMS voltage.
$ws.Range("C13").Formula = "=C12*C12/H3"
$ws.Range("C13").NumberFormat = "0.000"

# Re-enter E5's formula (value is unchanged) - in the authored edit this
# cell was individually retyped, breaking it out of the shared-formula
# group that otherwise spans E3:E9.
$ws.Range("E5").Formula = '=C5*C5/$H$3'

# --- Hidden chart-tracking range names left behind by Excel (macOS
#     Excel keeps these even without a live chart once a sheet has been
#     used as a chart source).
function Add-HiddenName($name, $ref) {
    $n = $wb.Names.Add($name, $ref)
    $n.Visible = $false
}
Add-HiddenName "_xlchart.v1.0"  '=Planilha1!$A$1'
Add-HiddenName "_xlchart.v1.1"  '=Planilha1!$A$2:$A$9'
Add-HiddenName "_xlchart.v1.2"  '=Planilha1!$B$1'
Add-HiddenName "_xlchart.v1.3"  '=Planilha1!$B$2:$B$9'
Add-HiddenName "_xlchart.v1.4"  '=Planilha1!$C$1'
Add-HiddenName "_xlchart.v1.5"  '=Planilha1!$C$2:$C$9'
Add-HiddenName "_xlchart.v2.10" '=Planilha1!$C$1'
Add-HiddenName "_xlchart.v2.11" '=Planilha1!$C$2:$C$9'
Add-HiddenName "_xlchart.v2.6"  '=Planilha1!$A$1'
Add-HiddenName "_xlchart.v2.7"  '=Planilha1!$A$2:$A$9'
Add-HiddenName "_xlchart.v2.8"  '=Planilha1!$B$1'
Add-HiddenName "_xlchart.v2.9"  '=Planilha1!$B$2:$B$9'

# Move the active selection to E11, matching where the author ended up.
$ws.Range("E11").Select()
